$p = $ppt.ActivePresentation

# --- Slide 1 (title slide): merge "Organisations" + " Präsentation" runs
#     into a single run reading "Organisations-Präsentation" (keep the
#     formatting / rPr of the second run, which has no err="1").
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$tr1 = $subtitle.TextFrame.TextRange
$firstRunLen = "Organisations".Length
# Clear the first run's text entirely so only the second run (without
# err="1") survives, then set its text to the merged, hyphenated phrase.
$run1a = $tr1.Characters(1, $firstRunLen)
$run1a.Text = ""
$remaining = $tr1.Text
$run1b = $tr1.Characters(1, $remaining.Length)
$run1b.Text = "Organisations-Präsentation"

# --- Slide 2 (Inhaltsverzeichnis): fix typo "Grobes Programmmodel" -> "Grobes Programmmodell"
$s2 = $p.Slides.Item(2)
$body2 = $s2.Shapes.Item(2)
$tr2 = $body2.TextFrame.TextRange
[void]$tr2.Replace("Grobes Programmmodel", "Grobes Programmmodell")

# --- Slide 3 (Das Spiel): add missing commas
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2)
$tr3 = $body3.TextFrame.TextRange
[void]$tr3.Replace("Verliert Länge wenn Item aufgehoben wird", "Verliert Länge, wenn Item aufgehoben wird")
[void]$tr3.Replace(" wenn altes aufgehoben wurde", ", wenn altes aufgehoben wurde")

# --- Delete slide 6 ("Untergruppen") entirely; slide 7/8 shift up.
$p.Slides.Item(6).Delete()
